$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 42588.471574074072
$ws.Range("A6").NumberFormat = "m/d/yy h:mm"

$ws.Range("B6").Value = "Named"
$ws.Range("C6").Value = 11854
$ws.Range("D6").Value = 5914
$ws.Range("E6").Value = 380
$ws.Range("F6").Value = 37
$ws.Range("G6").Value = 22
$ws.Range("H6").Value = 62
$ws.Range("I6").Value = 37
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
